$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (scheduled GitHub Actions data pull).
# Column D ("Price") holds numeric-looking text (e.g. "9.30", "0.100",
# "43.691.66") that must stay literal text, matching the original
# inline-string cells -- Excel would otherwise auto-coerce plain
# decimals into floating point numbers and drop trailing zeros. Forcing
# the Text number format for the duration of the write, then
# restoring the Normal cell style, preserves the text while leaving
# the cell style untouched (matches the diff: no style attribute change).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.691.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.282.17"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.25%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.92%  "

$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.920"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.629.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.281.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.794.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.93%  "

$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.97%  "

$ws.Range("E27").Value = "  +1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.32%  "

$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0928"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.28%  "

$ws.Range("E35").Value = "  +4.88%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0387"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.85%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("E39").Value = "  +5.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +44.98%  "

$ws.Range("E48").Value = "  +4.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "
